# feat: add 2022-Q3 data
#
# 1) Insert a new "2022-Q3" sheet right after the "总计" sheet (i.e. right
#    before "2021-Q4"), populated with the new quarter's fund-holding data.
#    We build it by duplicating the "2021-Q4" sheet (so all header text /
#    styles come along for free) and then overwriting its data rows.
# 2) Insert a new row at the top of the "总计" (totals) sheet summarising
#    the 2022-Q3 data, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$value) {
    # Force the cell to be stored as text (matches the source sheets, where
    # numeric-looking figures such as "6.79" are kept as strings), then
    # strip the format change back out so no stray style index is left on
    # the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Part 1: new "2022-Q3" worksheet
# ---------------------------------------------------------------------

$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Copy($q4Sheet)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q3"

# The template sheet has 5 data rows (rows 2-6); 2022-Q3 only needs 2, so
# drop rows 4-6 and overwrite rows 2-3 with the new figures.
$newSheet.Rows("4:6").Delete()

# Row 2 -> fund 002076
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "002076"
Set-TextValue $newSheet.Range("C2") "浙商中证500指数增强A"
Set-TextValue $newSheet.Range("D2") "6.79"
Set-TextValue $newSheet.Range("E2") "85.72"
Set-TextValue $newSheet.Range("F2") "1.11"
Set-TextValue $newSheet.Range("G2") "0.0754"
$newSheet.Range("H2").Value = 6

# Row 3 -> fund 007386
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "007386"
Set-TextValue $newSheet.Range("C3") "浙商中证500指数增强C"
Set-TextValue $newSheet.Range("D3") "1.97"
Set-TextValue $newSheet.Range("E3") "85.72"
Set-TextValue $newSheet.Range("F3") "1.11"
Set-TextValue $newSheet.Range("G3") "0.0219"
$newSheet.Range("H3").Value = 6

# Restore the originally-active sheet (the copy/rename above shifts focus
# onto the new sheet).
$wb.Worksheets.Item("2021-Q2").Activate()

# ---------------------------------------------------------------------
# Part 2: "总计" (totals) sheet gains a 2022-Q3 row at the top
# ---------------------------------------------------------------------

$total = $wb.Worksheets.Item("总计")

# Shift the 3 existing data rows down by one (copy bottom-up to preserve
# per-cell formatting along with the values).
$total.Range("A4:D4").Copy($total.Range("A5:D5"))
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.1

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
